$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '66.980.76'
$ws.Range('E2').Value = '  +0.09%  '
Set-TextValue $ws 'D3' '3.123.03'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue $ws 'D5' '577.98'
$ws.Range('E5').Value = '  -0.33%  '
Set-TextValue $ws 'D6' '173.21'
$ws.Range('E6').Value = '  +2.72%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.51%  '
Set-TextValue $ws 'D9' '6.45'
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('E12').Value = '  -0.78%  '
Set-TextValue $ws 'D13' '37.22'
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('E14').Value = '  -1.23%  '
Set-TextValue $ws 'D15' '3.642.28'
$ws.Range('E15').Value = '  +1.12%  '
Set-TextValue $ws 'D16' '66.964.43'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('E17').Value = '  -0.25%  '
Set-TextValue $ws 'D18' '3.123.87'
$ws.Range('E18').Value = '  +1.13%  '
Set-TextValue $ws 'D19' '16.28'
$ws.Range('E19').Value = '  +0.78%  '
Set-TextValue $ws 'D20' '476.96'
$ws.Range('E20').Value = '  +2.13%  '
Set-TextValue $ws 'D21' '0.710'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('E22').Value = '  +5.40%  '
Set-TextValue $ws 'D23' '83.98'
$ws.Range('E23').Value = '  +0.41%  '
Set-TextValue $ws 'D24' '13.31'
$ws.Range('E24').Value = '  +1.64%  '
Set-TextValue $ws 'D25' '2.29'
$ws.Range('E25').Value = '  -2.69%  '
Set-TextValue $ws 'D26' '10.14'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  +0.00%  '
Set-TextValue $ws 'D28' '7.92'
$ws.Range('E28').Value = '  -1.09%  '
Set-TextValue $ws 'D29' '2.37'
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('E30').Value = '  +0.19%  '
Set-TextValue $ws 'D31' '28.59'
$ws.Range('E31').Value = '  +1.18%  '
$ws.Range('E32').Value = '  +0.63%  '
Set-TextValue $ws 'D33' '0.0₃0955'
$ws.Range('E33').Value = '  -7.00%  '
$ws.Range('E34').Value = '  +0.09%  '
Set-TextValue $ws 'D35' '5.86'
$ws.Range('E35').Value = '  -0.58%  '
Set-TextValue $ws 'D36' '0.976'
$ws.Range('E36').Value = '  -2.73%  '
$ws.Range('E37').Value = '  +0.60%  '
Set-TextValue $ws 'D38' '2.06'
$ws.Range('E38').Value = '  -2.11%  '
Set-TextValue $ws 'D39' '50.20'
$ws.Range('E39').Value = '  -0.11%  '
Set-TextValue $ws 'D40' '0.312'
$ws.Range('E40').Value = '  -1.94%  '
$ws.Range('E41').Value = '  +0.98%  '
Set-TextValue $ws 'D42' '8.67'
$ws.Range('E42').Value = '  -0.12%  '
Set-TextValue $ws 'D43' '2.813.23'
$ws.Range('E43').Value = '  +1.36%  '
Set-TextValue $ws 'D44' '383.95'
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('E45').Value = '  -1.45%  '
Set-TextValue $ws 'D46' '2.57'
$ws.Range('E46').Value = '  -9.16%  '
Set-TextValue $ws 'D47' '135.72'
Set-TextValue $ws 'D49' '24.93'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('E51').Value = '  -0.60%  '
